$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 2335.8572
$ws.Range("I9").Value = 187.2
$ws.Range("K9").Value = 187.2
$ws.Range("M9").Value = -18.19999999999999
$ws.Range("H17").Value = 1789.8
$ws.Range("J17").Value = 1789.8
$ws.Range("L17").Value = 5369.4
$ws.Range("N17").Value = -5705.4
$ws.Range("H19").Value = 1482.6
$ws.Range("J19").Value = 1593.5555
$ws.Range("L19").Value = 1593.5555
$ws.Range("N19").Value = -1943.5555
$ws.Range("H42").Value = 1533
$ws.Range("I42").Value = 239.8
$ws.Range("J42").Value = 7999
$ws.Range("K42").Value = 719.4000000000001
$ws.Range("L42").Value = 23997
$ws.Range("M42").Value = -489.4000000000001
$ws.Range("N42").Value = -24457
$ws.Range("H43").Value = 3108.4443
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 3108.4443
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 3108.4443
$ws.Range("M43").ClearContents()
$ws.Range("N43").Value = -3246.4443
$ws.Range("H80").Value = 3925.1428
$ws.Range("I80").Value = 722.3333
$ws.Range("J80").Value = 4798.636
$ws.Range("K80").Value = 2166.9999
$ws.Range("L80").Value = 14395.908
$ws.Range("M80").Value = -1168.9999
$ws.Range("N80").Value = -16391.908
$ws.Range("H83").Value = 3925.1428
$ws.Range("I83").Value = 722.3333
$ws.Range("J83").Value = 4798.636
$ws.Range("K83").Value = 6500.9997
$ws.Range("L83").Value = 43187.724
$ws.Range("M83").Value = -1508.9997
$ws.Range("N83").Value = -53171.724
$ws.Range("H92").Value = 246.72728
$ws.Range("I92").Value = 258.35715
$ws.Range("J92").Value = 226.375
$ws.Range("K92").Value = 258.35715
$ws.Range("L92").Value = 226.375
$ws.Range("M92").Value = 989.64285
$ws.Range("N92").Value = -2722.375
$ws.Range("H96").Value = 1133.875
$ws.Range("I96").Value = 913
$ws.Range("J96").Value = 1354.75
$ws.Range("K96").Value = 2739
$ws.Range("L96").Value = 4064.25
$ws.Range("M96").Value = -1366
$ws.Range("N96").Value = -6810.25
$ws.Range("H100").Value = 2665.5
$ws.Range("I100").Value = 1832.625
$ws.Range("J100").Value = 5997
$ws.Range("K100").Value = 1832.625
$ws.Range("L100").Value = 5997
$ws.Range("M100").Value = -1291.625
$ws.Range("N100").Value = -7079
$ws.Range("H116").Value = 145772.14
$ws.Range("J116").Value = 3401.2
$ws.Range("L116").Value = 3401.2
$ws.Range("N116").Value = -10285.2
$ws.Range("H132").Value = 40241.31
$ws.Range("I132").Value = 41091
$ws.Range("K132").Value = 123273
$ws.Range("M132").Value = -120743
$ws.Range("H138").Value = 3199.56
$ws.Range("J138").Value = 3276.6667
$ws.Range("L138").Value = 9830.000100000001
$ws.Range("N138").Value = -20110.0001

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2013
$ws.Range("I2").Value = 1412.2307
$ws.Range("J2").Value = 2989.25
$ws.Range("K2").Value = 1412.2307
$ws.Range("L2").Value = 2989.25
$ws.Range("M2").Value = -1299.2307
$ws.Range("N2").Value = -3215.25
$ws.Range("H41").Value = 5291.4
$ws.Range("I41").Value = 5291.4
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 5291.4
$ws.Range("L41").Value = 0
$ws.Range("M41").Value = -4877.4
$ws.Range("N41").ClearContents()
$ws.Range("H110").Value = 2072.923
$ws.Range("I110").Value = 427.54544
$ws.Range("J110").Value = 11122.5
$ws.Range("K110").Value = 427.54544
$ws.Range("L110").Value = 11122.5
$ws.Range("M110").Value = 1617.45456
$ws.Range("N110").Value = -15212.5
$ws.Range("H116").Value = 2013
$ws.Range("I116").Value = 1412.2307
$ws.Range("J116").Value = 2989.25
$ws.Range("K116").Value = 1412.2307
$ws.Range("L116").Value = 2989.25
$ws.Range("M116").Value = 881.7692999999999
$ws.Range("N116").Value = -7577.25

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2013
$ws.Range("I3").Value = 1412.2307
$ws.Range("J3").Value = 2989.25
$ws.Range("K3").Value = 1412.2307
$ws.Range("L3").Value = 2989.25
$ws.Range("M3").Value = -1298.2307
$ws.Range("N3").Value = -3217.25
$ws.Range("H20").Value = 840.7222
$ws.Range("I20").Value = 879.6923
$ws.Range("K20").Value = 879.6923
$ws.Range("M20").Value = -632.6923

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1848
$ws.Range("I16").Value = 1468.5714
$ws.Range("K16").Value = 1468.5714
$ws.Range("M16").Value = -1181.5714
$ws.Range("H23").Value = 8556.286
$ws.Range("I23").Value = 5815.6665
$ws.Range("K23").Value = 5815.6665
$ws.Range("M23").Value = -5575.6665
$ws.Range("H27").Value = 8556.286
$ws.Range("I27").Value = 5815.6665
$ws.Range("K27").Value = 5815.6665
$ws.Range("M27").Value = -5623.6665
$ws.Range("H31").Value = 1289.4117
$ws.Range("I31").Value = 1325.7778
$ws.Range("J31").Value = 1248.5
$ws.Range("K31").Value = 1325.7778
$ws.Range("L31").Value = 1248.5
$ws.Range("M31").Value = -1030.7778
$ws.Range("N31").Value = -1838.5
$ws.Range("H34").Value = 1289.4117
$ws.Range("I34").Value = 1325.7778
$ws.Range("J34").Value = 1248.5
$ws.Range("K34").Value = 1325.7778
$ws.Range("L34").Value = 1248.5
$ws.Range("M34").Value = -1123.7778
$ws.Range("N34").Value = -1652.5
$ws.Range("H86").Value = 3999.6667
$ws.Range("I86").Value = 3999
$ws.Range("K86").Value = 3999
$ws.Range("M86").Value = -2876
$ws.Range("H89").Value = 3999.6667
$ws.Range("I89").Value = 3999
$ws.Range("K89").Value = 19995
$ws.Range("M89").Value = -14379
$ws.Range("H113").Value = 1848
$ws.Range("I113").Value = 1468.5714
$ws.Range("K113").Value = 1468.5714
$ws.Range("M113").Value = 701.4286
$ws.Range("H132").Value = 3901.8635
$ws.Range("I132").Value = 3759.5952
$ws.Range("K132").Value = 11278.7856
$ws.Range("M132").Value = -8748.785600000001

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H117").Value = 3545.25
$ws.Range("I117").Value = 900
$ws.Range("J117").Value = 3923.1428
$ws.Range("K117").Value = 2700
$ws.Range("L117").Value = 11769.4284
$ws.Range("M117").Value = 742
$ws.Range("N117").Value = -18653.4284
$ws.Range("H131").Value = 728565.25
$ws.Range("J131").Value = 1654488.6
$ws.Range("L131").Value = 4963465.800000001
$ws.Range("N131").Value = -4973545.800000001

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 10736.25
$ws.Range("I70").Value = 10675.615
$ws.Range("J70").Value = 10999
$ws.Range("K70").Value = 10675.615
$ws.Range("L70").Value = 10999
$ws.Range("M70").Value = -10405.615
$ws.Range("N70").Value = -11539
$ws.Range("H73").Value = 10736.25
$ws.Range("I73").Value = 10675.615
$ws.Range("J73").Value = 10999
$ws.Range("K73").Value = 10675.615
$ws.Range("L73").Value = 10999
$ws.Range("M73").Value = -9739.615
$ws.Range("N73").Value = -12871
$ws.Range("H80").Value = 7365.3335
$ws.Range("I80").Value = 3798.25
$ws.Range("J80").Value = 14499.5
$ws.Range("K80").Value = 3798.25
$ws.Range("L80").Value = 14499.5
$ws.Range("M80").Value = -2800.25
$ws.Range("N80").Value = -16495.5
$ws.Range("H83").Value = 7365.3335
$ws.Range("I83").Value = 3798.25
$ws.Range("J83").Value = 14499.5
$ws.Range("K83").Value = 18991.25
$ws.Range("L83").Value = 72497.5
$ws.Range("M83").Value = -13999.25
$ws.Range("N83").Value = -82481.5

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3617.7144
$ws.Range("I40").Value = 3554.25
$ws.Range("K40").Value = 3554.25
$ws.Range("M40").Value = -3418.25
$ws.Range("H122").Value = 3428.2856
$ws.Range("I122").Value = 3128.4
$ws.Range("K122").Value = 9385.200000000001
$ws.Range("M122").Value = -6935.200000000001

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").ClearContents()
$ws.Range("H132").Value = 974.9231
$ws.Range("I132").Value = 974.9231
$ws.Range("K132").Value = 2924.7693
$ws.Range("M132").Value = -394.7692999999999
$ws.Range("H136").Value = 2249.3333
$ws.Range("J136").Value = 2664.3333
$ws.Range("L136").Value = 7992.999899999999
$ws.Range("N136").Value = -13092.9999
